$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row renames: the merged-AHB diff columns were regenerated from
# "_old" / "_new" suffixes to the concrete version tags "_FV2404" / "_FV2410".
$headerRenames = @{
    "A1" = "Segmentname_FV2404";
    "B1" = "Segmentgruppe_FV2404";
    "C1" = "Segment_FV2404";
    "D1" = "Datenelement_FV2404";
    "E1" = "Segment ID_FV2404";
    "F1" = "Code_FV2404";
    "G1" = "Qualifier_FV2404";
    "H1" = "Beschreibung_FV2404";
    "I1" = "Bedingungsausdruck_FV2404";
    "J1" = "Bedingung_FV2404";
    "L1" = "Segmentname_FV2410";
    "M1" = "Segmentgruppe_FV2410";
    "N1" = "Segment_FV2410";
    "O1" = "Datenelement_FV2410";
    "P1" = "Segment ID_FV2410";
    "Q1" = "Code_FV2410";
    "R1" = "Qualifier_FV2410";
    "S1" = "Beschreibung_FV2410";
    "T1" = "Bedingungsausdruck_FV2410";
    "U1" = "Bedingung_FV2410";
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value2 = $headerRenames[$addr]
}

# Turn the data range into a table (ListObject)
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U57"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
